# Add a second worksheet ("Sheet2") positioned after Sheet1, populate it with
# the two test strings, and restore Sheet1 as the active/selected sheet
# (matching the saved selection state: Sheet2!F5 remembered, Sheet1 active).

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet2"

$ws2.Range("A1").Value = "sheet_test"
$ws2.Range("A2").Value = "this is just a test for passing Excel sheets as a readit argument"

# Leave Sheet2's remembered selection at F5, then reactivate Sheet1 so it
# stays the tab Excel shows on open.
$ws2.Range("F5").Select()
$ws1.Select()
